$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 151-152, shifting the existing rows 151-249 down to 153-251
$ws.Rows("151:152").Insert()

# Row 151
$ws.Range("A151").Value = 7
$ws.Range("B151").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C151").Value = 'Ñuble'
$ws.Range("D151").Value = 44767
$ws.Range("E151").Value = 16
$ws.Range("F151").Value = 100112009
$ws.Range("G151").Value = 'Acelga'
$ws.Range("H151").Value = 'Sin especificar'
$ws.Range("I151").Value = 'Primera'
$ws.Range("J151").Value = 200
$ws.Range("K151").Value = 700
$ws.Range("L151").Value = 800
$ws.Range("M151").Value = 750
$ws.Range("N151").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O151").Value = 'Provincia de Diguillín'
$ws.Range("P151").Value = 750
$ws.Range("Q151").Value = 1
$ws.Range("R151").Value = 'Hortaliza'

# Row 152
$ws.Range("A152").Value = 7
$ws.Range("B152").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C152").Value = 'Ñuble'
$ws.Range("D152").Value = 44767
$ws.Range("E152").Value = 16
$ws.Range("F152").Value = 100112009
$ws.Range("G152").Value = 'Acelga'
$ws.Range("H152").Value = 'Sin especificar'
$ws.Range("I152").Value = 'Segunda'
$ws.Range("J152").Value = 120
$ws.Range("K152").Value = 600
$ws.Range("L152").Value = 600
$ws.Range("M152").Value = 600
$ws.Range("N152").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O152").Value = 'Provincia de Diguillín'
$ws.Range("P152").Value = 600
$ws.Range("Q152").Value = 1
$ws.Range("R152").Value = 'Hortaliza'
